# Update "想去人数" (want-to-go headcount) figures on the refreshed
# scrape output for both the "展览" sheet and its mirror in "全部类型".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "全部类型") {
        $ws.Range("F2").Value = 661
        $ws.Range("F3").Value = 503
        $ws.Range("F7").Value = 42
        $ws.Range("F8").Value = 2440
        $ws.Range("F9").Value = 4137
        $ws.Range("F10").Value = 99
    } else {
        $ws.Range("F2").Value = 661
        $ws.Range("F3").Value = 503
        $ws.Range("F7").Value = 42
        $ws.Range("F8").Value = 2439
        $ws.Range("F9").Value = 4137
        $ws.Range("F10").Value = 99
    }
}
